$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "update version to v0.4 structure: target=common options, value=options
# by case" - the two syntax-highlighted JSON "assert" payloads in column C swap
# their key names:
#   C2 was {"key":"age"}    -> becomes {"target":"age"}
#   C3 was {"value":"44"}   -> stays   {"value":"44"}  (same literal text)
#
# C3's displayed text doesn't actually change, so it is left alone. C2's first
# field name changes from "key" to "target"; we rebuild it with Characters()
# runs so the JSON-highlighter rich-text coloring (purple keys / dark
# punctuation / green values, with the mixed "Sarasa Mono CL" + "ＭＳ Ｐゴシック"
# font pairing already used by the sheet) is preserved exactly like the other
# untouched cell.

$c2 = $ws.Range("C2")
$c2.Value2 = '{"target":"age"}'

$purple = 9703559   # RGB FF871094 (key color)
$dark   = 526344    # RGB FF080808 (punctuation color)
$green  = 1539334   # RGB FF067D17 (value color)

function Set-RunStyle($cell, $start, $len, $color, $fontName) {
    $run = $cell.Characters($start, $len)
    $run.Font.Size = 11
    $run.Font.Color = $color
    $run.Font.Name = $fontName
}

# {"target":"age"}
#  1  2345678 9  10 11  121314  15 16
#  {  "target"  "   :   "  age   "   }
Set-RunStyle $c2 2  1 $purple "Sarasa Mono CL"     # opening quote of key
Set-RunStyle $c2 3  6 $purple "ＭＳ Ｐゴシック"      # target
Set-RunStyle $c2 9  1 $purple "Sarasa Mono CL"     # closing quote of key
Set-RunStyle $c2 10 1 $dark   "Sarasa Mono CL"     # :
Set-RunStyle $c2 11 1 $green  "Sarasa Mono CL"     # opening quote of value
Set-RunStyle $c2 12 3 $green  "ＭＳ Ｐゴシック"      # age
Set-RunStyle $c2 15 1 $green  "Sarasa Mono CL"     # closing quote of value
Set-RunStyle $c2 16 1 $dark   "Sarasa Mono CL"     # }

# Move the active selection, matching the saved workbook's cursor position.
$ws.Range("G8").Select()
